$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (values scraped on 2023-01-16)
$updates = @{
    "D2" = "300.10"
    "E2" = "-0.54%"
    "D3" = "31.64"
    "E3" = "0.65%"
    "D4" = "5.137"
    "E4" = "0.27%"
    "D5" = "0.08113"
    "E5" = "10.41%"
    "D6" = "2.620"
    "E6" = "23.46%"
    "D7" = "7.783"
    "E7" = "-1.84%"
    "D8" = "3.906"
    "E8" = "2.07%"
    "D9" = "0.9314"
    "E9" = "1.28%"
    "D10" = "0.1761"
    "E10" = "3.20%"
    "D11" = "0.07357"
    "E11" = "-2.01%"
    "D12" = "0.08828"
    "E12" = "8.25%"
    "D13" = "0.03027"
    "E13" = "0.05%"
    "D14" = "0.1001"
    "E14" = "0.60%"
    "D15" = "0.001508"
    "E15" = "-0.27%"
    "D16" = "0.005878"
    "E16" = "-4.33%"
    "D17" = "3.570"
    "E17" = "3.43%"
    "D18" = "2.286"
    "E18" = "2.90%"
    "D19" = "0.3272"
    "E19" = "-0.29%"
    "E20" = "-0.84%"
    "D21" = "4.166"
    "E21" = "-10.46%"
    "D22" = "0.1681"
    "E22" = "7.17%"
    "D23" = "0.04615"
    "E23" = "-0.47%"
    "D24" = "0.001239"
    "E24" = "0.80%"
    "D25" = "0.004531"
    "E25" = "1.41%"
    "E26" = "-7.63%"
    "D27" = "0.0003410"
    "E27" = "-0.79%"
    "D39" = "0.01768"
    "E39" = "2.29%"
    "D40" = "0.04607"
    "E40" = "1.96%"
    "D41" = "0.006942"
    "E41" = "-4.50%"
    "E42" = "2.30%"
    "D43" = "0.002190"
    "E43" = "-1.72%"
    "D44" = "0.01037"
    "E44" = "-2.85%"
    "D45" = "0.00006293"
    "E45" = "0.00%"
    "E46" = "0.01%"
    "D47" = "0.008400"
    "E47" = "-15.93%"
    "D48" = "0.7484"
    "E48" = "-7.43%"
    "D49" = "0.00002100"
    "E49" = "0.01%"
    "D50" = "0.0002000"
    "E50" = "0.08%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = $origStyle
}
